# Admin catalog product update - Testdata update
#
# 1. Rename "Sheet1" -> "Catalog"
# 2. Move the active/selected tab from "Cardtile,PLP,CLP" to "Catalog"
# 3. Update the selection on the newly active "Catalog" sheet to DI2
#    (was D5), leaving the selection on "Cardtile,PLP,CLP" untouched (AP5)

$wb = $excel.ActiveWorkbook

# --- Rename Sheet1 -> Catalog ------------------------------------------------
$wsCatalog = $wb.Worksheets.Item("Sheet1")
$wsCatalog.Name = "Catalog"

# --- Move the active tab to Catalog and update its selection ---------------
$wsCatalog.Activate()
$wsCatalog.Range("DI2").Select()
